$d = $word.ActiveDocument

# Locate the first "hola" (the Heading1 question title at the top of the doc)
# and turn it into "Hola?" -- but keep it split across three runs (H | ola | ?)
# that all resolve to the same run properties, matching the target markup.
$found = $d.Content
$found.Find.Execute("hola", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$start = $found.Start
$end = $found.End

# Insert the trailing "?" right after "hola" first, while the whole run is
# still contiguous, then give the new character a throwaway distinguishing
# format so it doesn't get silently re-merged into its neighbour.
$qMark = $d.Range($end, $end)
$qMark.InsertBefore("?")
$qRange = $d.Range($end, $end + 1)
$qRange.Bold = 1

# Capitalize the leading "h" -> "H" in its own run, using the same
# throwaway-format trick to keep it distinct from the rest of the word.
$hRange = $d.Range($start, $start + 1)
$hRange.Bold = 1
$hRange.Text = "H"

# Give the middle "ola" the same throwaway format so all three pieces are
# separate runs once the format gets normalized back below.
$olaRange = $d.Range($start + 1, $end)
$olaRange.Bold = 1

# Drop the throwaway formatting so every run ends up with the original
# (sz/szCs-only) run properties, while remaining three distinct runs.
$wholeRange = $d.Range($start, $end + 1)
$wholeRange.Bold = 0
